$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 144, shifting existing rows 144:230 down to 145:231
$ws.Rows("144").Insert()

# Populate the newly inserted row 144 with its data
$ws.Range("A144").Value = 10
$ws.Range("B144").Value = "Vega Modelo de Temuco"
$ws.Range("C144").Value = "La Araucanía"
$ws.Range("D144").Value = 44596
$ws.Range("E144").Value = 9
$ws.Range("F144").Value = 100112001
$ws.Range("G144").Value = "Berenjena"
$ws.Range("H144").Value = "Sin especificar"
$ws.Range("I144").Value = "Primera"
$ws.Range("J144").Value = 70
$ws.Range("K144").Value = 10000
$ws.Range("L144").Value = 12000
$ws.Range("M144").Value = 11143
$ws.Range("N144").Value = "`$/caja 60 unidades"
$ws.Range("O144").Value = "Región del Maule"
$ws.Range("P144").Value = 186
$ws.Range("Q144").Value = 60
$ws.Range("R144").Value = "Hortaliza"
